$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts -> K) values for column G, rows 2-15
$values = @{
    2  = 2
    3  = 4
    4  = 0
    5  = 8
    6  = 3
    7  = 1
    8  = 2
    9  = 3
    10 = 3
    11 = 2
    12 = 2
    13 = 2
    14 = 3
    15 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
